$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.946.68"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.877.73"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4923"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2921"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06635"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").Value = "1.878.44"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07236"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6659"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.880"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "29.929.38"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007880"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9992"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").Value = "2.119.15"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9986"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.779"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.787"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.052"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -4.29%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.190"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08742"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.971"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05045"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7126"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.114"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01771"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.689"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.180"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9306"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4243"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.761"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.67%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9987"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("E45").Value = "  -2.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1269"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05662"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3786"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.280"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.17%  "
